# income_statement_dashboard_data.xlsx — v1 dashboard edit
# Mirrors the author's edit: a few actuals inputs were revised, a spacer row
# was inserted above "Gross Profit" (pushing everything below it down one
# row), and the "actuals" tab became the active/selected sheet (cell E3).

$wb = $excel.ActiveWorkbook

$actuals = $wb.Worksheets.Item("actuals")
$targets = $wb.Worksheets.Item("targets")

# ---------------------------------------------------------------------
# 1. Revise a handful of input cells on "actuals" (pre-insert row numbers).
#    Everything else on the sheet (Total Revenue, Total COGS, Gross
#    Profit, margins, Net Operating Income, Net Income, ...) recalculates
#    automatically from these.
# ---------------------------------------------------------------------
$actuals.Range("C2").Value2 = 3937368.5
$actuals.Range("D2").Value2 = 3523797.5

$actuals.Range("C7").Value2 = 1600200
$actuals.Range("D7").Value2 = 1500640
$actuals.Range("E7").Value2 = 1420027.5
$actuals.Range("F7").Value2 = 1388595

$actuals.Range("F8").Value2 = 305319
$actuals.Range("F9").Value2 = 201546

# ---------------------------------------------------------------------
# 2. Insert a new blank spacer row above "Gross Profit" (old row 11),
#    shifting Gross Profit and every row below it down by one.
# ---------------------------------------------------------------------
$actuals.Rows("11:11").Insert()

# Re-assert the "Average" column formula across the now-contiguous
# P4Y-average rows (Cost of Product Sales through Gross Profit Margin),
# skipping the new blank spacer row, matching the author's re-fill.
$actuals.Range("B7:B13").Formula = "=AVERAGE(C7:F7)"
$actuals.Range("B11").ClearContents()

# ---------------------------------------------------------------------
# 3. "targets": re-enter the decay formulas across D2:F4 as one fill so
#    they form a single contiguous formula block (values are unchanged).
# ---------------------------------------------------------------------
$targets.Range("D2:F4").Formula = "=C2*0.95"

# ---------------------------------------------------------------------
# 4. Make "actuals" the active sheet / active tab, with E3 selected
#    (previously "targets" was active with D20 selected on "actuals").
# ---------------------------------------------------------------------
$actuals.Activate()
$actuals.Range("E3").Select()
